$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header updates
$ws.Range("A8").Value = "Volume 32   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/29/2025  Through  10/5/2025"

# Row 14
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = -50
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = -16.666666666666
$ws.Range("I14").Value = 28
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = -6.666666666666
$ws.Range("L14").Value = 115.384615384615
$ws.Range("M14").Value = -49.090909090909
$ws.Range("N14").Value = -75.438596491228

# Row 15
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -42.857142857142
$ws.Range("F15").Value = 17
$ws.Range("G15").Value = 29
$ws.Range("H15").Value = -41.379310344827
$ws.Range("I15").Value = 172
$ws.Range("J15").Value = 153
$ws.Range("K15").Value = 12.418300653594
$ws.Range("L15").Value = 50.877192982456
$ws.Range("M15").Value = 72
$ws.Range("N15").Value = -37.681159420289

# Row 16
$ws.Range("C16").Value = 28
$ws.Range("D16").Value = 33
$ws.Range("E16").Value = -15.151515151515
$ws.Range("F16").Value = 112
$ws.Range("G16").Value = 107
$ws.Range("H16").Value = 4.672897196261
$ws.Range("I16").Value = 921
$ws.Range("J16").Value = 1082
$ws.Range("K16").Value = -14.879852125693
$ws.Range("L16").Value = -11.612284069097
$ws.Range("M16").Value = -42.652552926525
$ws.Range("N16").Value = -84.310051107325

# Row 17
$ws.Range("C17").Value = 60
$ws.Range("D17").Value = 41
$ws.Range("E17").Value = 46.341463414634
$ws.Range("G17").Value = 256
$ws.Range("H17").Value = -12.890625
$ws.Range("I17").Value = 2296
$ws.Range("J17").Value = 2436
$ws.Range("K17").Value = -5.747126436781
$ws.Range("L17").Value = 4.936014625228
$ws.Range("M17").Value = 77.16049382716
$ws.Range("N17").Value = -20.663441603317

# Row 18
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 22
$ws.Range("E18").Value = -68.181818181818
$ws.Range("F18").Value = 64
$ws.Range("G18").Value = 100
$ws.Range("H18").Value = -36
$ws.Range("I18").Value = 781
$ws.Range("J18").Value = 825
$ws.Range("K18").Value = -5.333333333333
$ws.Range("L18").Value = -10.538373424971
$ws.Range("M18").Value = -50.942211055276
$ws.Range("N18").Value = -88.814093382984

# Row 19
$ws.Range("C19").Value = 89
$ws.Range("D19").Value = 48
$ws.Range("E19").Value = 85.416666666666
$ws.Range("F19").Value = 332
$ws.Range("G19").Value = 241
$ws.Range("H19").Value = 37.759336099585
$ws.Range("I19").Value = 2599
$ws.Range("J19").Value = 2529
$ws.Range("K19").Value = 2.767892447607
$ws.Range("L19").Value = -7.377049180327
$ws.Range("M19").Value = 22.363465160075
$ws.Range("N19").Value = -58.954516740366

# Row 20
$ws.Range("C20").Value = 39
$ws.Range("D20").Value = 49
$ws.Range("E20").Value = -20.408163265306
$ws.Range("F20").Value = 151
$ws.Range("G20").Value = 177
$ws.Range("H20").Value = -14.689265536723
$ws.Range("I20").Value = 1467
$ws.Range("J20").Value = 1515
$ws.Range("K20").Value = -3.168316831683
$ws.Range("L20").Value = 5.236728837876
$ws.Range("M20").Value = 8.666666666666
$ws.Range("N20").Value = -89.931365820178

# Row 21
$ws.Range("C21").Value = 228
$ws.Range("D21").Value = 202
$ws.Range("E21").Value = 12.871287128712
$ws.Range("F21").Value = 904
$ws.Range("G21").Value = 916
$ws.Range("H21").Value = -1.310043668122
$ws.Range("I21").Value = 8264
$ws.Range("J21").Value = 8570
$ws.Range("K21").Value = -3.570595099183
$ws.Range("L21").Value = -1.969157769869
$ws.Range("M21").Value = 1.735811892158
$ws.Range("N21").Value = -77.687780117716

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = -35.714285714285
$ws.Range("I22").Value = 89
$ws.Range("J22").Value = 104
$ws.Range("K22").Value = -14.423076923076
$ws.Range("L22").Value = -6.315789473684
$ws.Range("M22").Value = -1.111111111111

# Row 23
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = -14.285714285714
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 21
$ws.Range("H23").Value = -42.857142857142
$ws.Range("I23").Value = 173
$ws.Range("J23").Value = 178
$ws.Range("K23").Value = -2.808988764044
$ws.Range("L23").Value = -4.945054945054
$ws.Range("M23").Value = 50.434782608695

# Row 24
$ws.Range("C24").Value = 169
$ws.Range("D24").Value = 210
$ws.Range("E24").Value = -19.523809523809
$ws.Range("F24").Value = 764
$ws.Range("G24").Value = 834
$ws.Range("H24").Value = -8.393285371702
$ws.Range("I24").Value = 7070
$ws.Range("J24").Value = 7164
$ws.Range("K24").Value = -1.312116136236
$ws.Range("L24").Value = 0.028296547821
$ws.Range("M24").Value = 49.978786593126

# Row 25
$ws.Range("C25").Value = 65
$ws.Range("D25").Value = 102
$ws.Range("E25").Value = -36.274509803921
$ws.Range("F25").Value = 319
$ws.Range("G25").Value = 393
$ws.Range("H25").Value = -18.82951653944
$ws.Range("I25").Value = 2818
$ws.Range("J25").Value = 3107
$ws.Range("K25").Value = -9.301577084003
$ws.Range("L25").Value = 14.739413680781

# Row 26
$ws.Range("C26").Value = 102
$ws.Range("D26").Value = 92
$ws.Range("E26").Value = 10.869565217391
$ws.Range("F26").Value = 391
$ws.Range("G26").Value = 449
$ws.Range("H26").Value = -12.917594654788
$ws.Range("I26").Value = 3798
$ws.Range("J26").Value = 4000
$ws.Range("K26").Value = -5.05
$ws.Range("L26").Value = 9.137931034482
$ws.Range("M26").Value = 8.483290488431

# Row 27
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 9
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Value = 22
$ws.Range("G27").Value = 38
$ws.Range("H27").Value = -42.105263157894
$ws.Range("I27").Value = 213
$ws.Range("J27").Value = 240
$ws.Range("K27").Value = -11.25
$ws.Range("L27").Value = 5.970149253731

# Row 28
$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 40
$ws.Range("F28").Value = 26
$ws.Range("G28").Value = 34
$ws.Range("H28").Value = -23.529411764705
$ws.Range("I28").Value = 339
$ws.Range("J28").Value = 338
$ws.Range("K28").Value = 0.295857988165
$ws.Range("L28").Value = 9.708737864077

# Row 29
$ws.Range("C29").Value = 2
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -50
$ws.Range("G29").Value = 20
$ws.Range("H29").Value = -80
$ws.Range("I29").Value = 57
$ws.Range("J29").Value = 93
$ws.Range("K29").Value = -38.709677419354
$ws.Range("L29").Value = -21.917808219178
$ws.Range("M29").Value = -62.745098039215
$ws.Range("N29").Value = -85.309278350515

# Row 30
$ws.Range("C30").Value = 2
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = -33.333333333333
$ws.Range("G30").Value = 15
$ws.Range("H30").Value = -73.333333333333
$ws.Range("I30").Value = 45
$ws.Range("J30").Value = 75
$ws.Range("K30").Value = -40
$ws.Range("L30").Value = -18.181818181818
$ws.Range("M30").Value = -62.809917355371
$ws.Range("N30").Value = -87.215909090909

# Row 31
$ws.Range("D31").Value = 2
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G31").Value = 6
$ws.Range("H31").Value = -83.333333333333
$ws.Range("I31").Value = 26
$ws.Range("J31").Value = 31
$ws.Range("K31").Value = -16.129032258064
$ws.Range("L31").Value = -39.53488372093

# Row 33
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "0"
$ws.Range("C33").NumberFormat = "General"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "***.*"
$ws.Range("E33").NumberFormat = "General"
$ws.Range("F33").Value = 5
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 400
$ws.Range("I33").Value = 25
$ws.Range("K33").Value = 31.578947368421
$ws.Range("L33").Value = -13.793103448275
